$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 291 (shifts old rows 291.. down to 293..)
$ws.Rows.Item(291).Resize(2).Insert()

# New row 291 values
$ws.Cells.Item(291, 1).Value = 6
$ws.Cells.Item(291, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(291, 3).Value = "Metropolitana"
$ws.Cells.Item(291, 4).Value = 44516
$ws.Cells.Item(291, 5).Value = 13
$ws.Cells.Item(291, 6).Value = 100112043
$ws.Cells.Item(291, 7).Value = "Pepino ensalada"
$ws.Cells.Item(291, 8).Value = "Sin especificar"
$ws.Cells.Item(291, 9).Value = "Primera"
$ws.Cells.Item(291, 10).Value = 1400
$ws.Cells.Item(291, 11).Value = 5000
$ws.Cells.Item(291, 12).Value = 6000
$ws.Cells.Item(291, 13).Value = 5464
$ws.Cells.Item(291, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(291, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(291, 16).Value = 91
$ws.Cells.Item(291, 17).Value = 60
$ws.Cells.Item(291, 18).Value = "Hortaliza"

# New row 292 values
$ws.Cells.Item(292, 1).Value = 6
$ws.Cells.Item(292, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(292, 3).Value = "Metropolitana"
$ws.Cells.Item(292, 4).Value = 44516
$ws.Cells.Item(292, 5).Value = 13
$ws.Cells.Item(292, 6).Value = 100112043
$ws.Cells.Item(292, 7).Value = "Pepino ensalada"
$ws.Cells.Item(292, 8).Value = "Sin especificar"
$ws.Cells.Item(292, 9).Value = "Segunda"
$ws.Cells.Item(292, 10).Value = 490
$ws.Cells.Item(292, 11).Value = 3500
$ws.Cells.Item(292, 12).Value = 4000
$ws.Cells.Item(292, 13).Value = 3827
$ws.Cells.Item(292, 14).Value = "$/caja 100 unidades"
$ws.Cells.Item(292, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(292, 16).Value = 38
$ws.Cells.Item(292, 17).Value = 100
$ws.Cells.Item(292, 18).Value = "Hortaliza"
